# Applies the Oct 3 2024 cryptos.xlsx price/volume refresh to the active worksheet.
# Column D prices that look like plain numbers are written via a text-formatted
# round-trip (NumberFormat "@") so Excel keeps them as literal strings (e.g. "318.80")
# instead of silently coercing them into numeric values; the cells original style
# is restored immediately afterwards so no visual formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

# Row 2
$ws.Range("D2").Value = '60.778.43'
$ws.Range("E2").Value = '  -0.13%  '
# Row 3
$ws.Range("D3").Value = '2.348.43'
$ws.Range("E3").Value = '  -1.20%  '
# Row 4
$ws.Range("E4").Value = '  -0.03%  '
# Row 5
Set-TextValue $ws.Range("D5") '544.02'
$ws.Range("E5").Value = '  +0.01%  '
# Row 6
Set-TextValue $ws.Range("D6") '137.02'
$ws.Range("E6").Value = '  -2.92%  '
# Row 7
$ws.Range("E7").Value = '  +0.01%  '
# Row 8
$ws.Range("E8").Value = '  -4.61%  '
# Row 9
$ws.Range("D9").Value = '2.346.39'
$ws.Range("E9").Value = '  -1.31%  '
# Row 10
$ws.Range("E10").Value = '  -0.28%  '
# Row 11
$ws.Range("E11").Value = '  +1.89%  '
# Row 12
$ws.Range("E12").Value = '  -0.63%  '
# Row 13
$ws.Range("E13").Value = '  -0.36%  '
# Row 14
$ws.Range("E14").Value = '  -3.01%  '
# Row 15
$ws.Range("D15").Value = '2.772.50'
$ws.Range("E15").Value = '  -1.15%  '
# Row 16
$ws.Range("D16").Value = '60.678.26'
$ws.Range("E16").Value = '  -0.18%  '
# Row 17
$ws.Range("E17").Value = '  -2.12%  '
# Row 18
$ws.Range("D18").Value = '2.347.67'
$ws.Range("E18").Value = '  -1.08%  '
# Row 19
$ws.Range("E19").Value = '  +0.27%  '
# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D20") '318.80'
$ws.Range("E20").Value = '  +0.71%  '
# Row 21
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D21") '4.12'
$ws.Range("E21").Value = '  +0.29%  '
# Row 22
$ws.Range("E22").Value = '  -2.33%  '
# Row 23
$ws.Range("E23").Value = '  -0.09%  '
# Row 24
Set-TextValue $ws.Range("D24") '63.33'
$ws.Range("E24").Value = '  +0.54%  '
# Row 25
Set-TextValue $ws.Range("D25") '1.68'
$ws.Range("E25").Value = '  -7.12%  '
# Row 26
Set-TextValue $ws.Range("D26") '8.32'
$ws.Range("E26").Value = '  +7.39%  '
# Row 27
$ws.Range("E27").Value = '  +0.03%  '
# Row 28
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '2.464.29'
$ws.Range("E28").Value = '  -1.11%  '
# Row 29
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D29") '7.96'
$ws.Range("E29").Value = '  -0.38%  '
# Row 30
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D30") '497.81'
$ws.Range("E30").Value = '  -3.72%  '
# Row 31
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D31") '1.37'
$ws.Range("E31").Value = '  -3.99%  '
# Row 32
$ws.Range("D32").Value = '0.0₃0860'
$ws.Range("E32").Value = '  -6.98%  '
# Row 33
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D33") '0.145'
$ws.Range("E33").Value = '  +0.59%  '
# Row 34
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D34") '1.79'
$ws.Range("E34").Value = '  -2.09%  '
# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D35") '1.50'
$ws.Range("E35").Value = '  -3.83%  '
# Row 36
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D36") '0.999'
$ws.Range("E36").Value = '  -0.02%  '
# Row 37
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D37") '4.58'
$ws.Range("E37").Value = '  -1.30%  '
# Row 38
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range("D38") '0.375'
$ws.Range("E38").Value = '  +0.19%  '
# Row 39
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D39") '18.45'
$ws.Range("E39").Value = '  +2.14%  '
# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D40") '1.83'
$ws.Range("E40").Value = '  +6.39%  '
# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range("D41") '5.24'
$ws.Range("E41").Value = '  -4.14%  '
# Row 42
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D42") '144.11'
$ws.Range("E42").Value = '  +5.57%  '
# Row 43
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D43") '0.999'
$ws.Range("E43").Value = '  -0.12%  '
# Row 44
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D44") '40.61'
$ws.Range("E44").Value = '  +0.90%  '
# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D45") '143.36'
$ws.Range("E45").Value = '  +2.92%  '
# Row 46
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D46") '3.56'
$ws.Range("E46").Value = '  +0.43%  '
# Row 47
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D47") '2.04'
$ws.Range("E47").Value = '  -8.95%  '
# Row 48
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D48") '0.0517'
$ws.Range("E48").Value = '  +0.49%  '
# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D49") '19.06'
$ws.Range("E49").Value = '  -6.31%  '
# Row 50
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D50") '0.568'
$ws.Range("E50").Value = '  -1.42%  '
# Row 51
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D51") '0.0901'
$ws.Range("E51").Value = '  -1.29%  '

Write-Host "Applied cryptos update"